# "Apply corrections to slides supplied by Mum"
#
# Fixes the lyric line on the "O sovereign God" song slide (slide 4):
#   "Encouraged by your spirit we may live," -> "...your Spirit we may live,"
# i.e. capitalise "spirit" -> "Spirit" (referring to the Holy Spirit), which
# splits the original single run into three runs (before / "Spirit " / after).

$p = $ppt.ActivePresentation

# Slide 4 ("4. Encouraged by your spirit we may live, ...") — Content Placeholder 1
$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

$fullText = $textRange.Text
$searchTerm = "spirit "
$startIndex = $fullText.IndexOf($searchTerm)

if ($startIndex -ge 0) {
    $wordRange = $textRange.Characters($startIndex + 1, $searchTerm.Length)
    $wordRange.Text = "Spirit "
}
